# Fruta / hortaliza, semanal
# Insert two new weekly data rows (374 and 375) into the Alcachofa sheet,
# pushing the existing rows 374-392 down to 376-394.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 374:375 (shifts old rows 374-392 down to 376-394)
$ws.Range("A374:A375").EntireRow.Insert()

# New row 374
$ws.Range("A374").Value = 9
$ws.Range("B374").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C374").Value = "Metropolitana"
$ws.Range("D374").Value = 44714
$ws.Range("E374").Value = 13
$ws.Range("F374").Value = 100112013
$ws.Range("G374").Value = "Alcachofa"
$ws.Range("H374").Value = "Española"
$ws.Range("I374").Value = "Primera"
$ws.Range("J374").Value = 52
$ws.Range("K374").Value = 24000
$ws.Range("L374").Value = 24000
$ws.Range("M374").Value = 24000
$ws.Range("N374").Value = "$/caja 30 unidades"
$ws.Range("O374").Value = "Provincia del Elquí"
$ws.Range("P374").Value = 800
$ws.Range("Q374").Value = 30
$ws.Range("R374").Value = "Hortaliza"

# New row 375
$ws.Range("A375").Value = 9
$ws.Range("B375").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C375").Value = "Metropolitana"
$ws.Range("D375").Value = 44714
$ws.Range("E375").Value = 13
$ws.Range("F375").Value = 100112013
$ws.Range("G375").Value = "Alcachofa"
$ws.Range("H375").Value = "Española"
$ws.Range("I375").Value = "Segunda"
$ws.Range("J375").Value = 25
$ws.Range("K375").Value = 22000
$ws.Range("L375").Value = 22000
$ws.Range("M375").Value = 22000
$ws.Range("N375").Value = "$/caja 40 unidades"
$ws.Range("O375").Value = "Provincia del Elquí"
$ws.Range("P375").Value = 550
$ws.Range("Q375").Value = 40
$ws.Range("R375").Value = "Hortaliza"
